# POM Updated script ( Create,delete ==> Cust. & project)
#
# Adds a new "taskspagedetails" worksheet (after the existing sheets)
# containing a small 2x2 table: headers "cutomerName" / "projectName"
# in row 1 and data "Microsoft" / "WebApplication" in row 2, formats it
# with a thin box border around every cell and a yellow fill on the
# header row, and makes it the active/selected sheet (mirroring the
# tabSelected move away from "managercreds").

$wb = $excel.ActiveWorkbook

# --- add the new sheet at the very end of the tab strip -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "taskspagedetails"

# --- header row -------------------------------------------------------------
$ws.Range("A1").Value = "cutomerName"
$ws.Range("B1").Value = "projectName"

# --- data row ----------------------------------------------------------------
$ws.Range("A2").Value = "Microsoft"
$ws.Range("B2").Value = "WebApplication"

# --- column widths (approximate the authored widths) ------------------------
$ws.Columns.Item(1).ColumnWidth = 11.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.333333333333334

# --- borders: thin box around every cell in the 2x2 table -------------------
$fullRange = $ws.Range("A1:B2")
$fullRange.Borders.LineStyle = 1
$fullRange.Borders.Weight = 2

# --- header fill: yellow ------------------------------------------------------
$headerRange = $ws.Range("A1:B1")
$headerRange.Interior.Color = 65535

# --- selection / active sheet -------------------------------------------------
[void]$ws.Activate()
[void]$ws.Range("G27").Select()
